$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Error Rate (column D) values for the affected rows, as reviewed.
$ws.Range("D3").Value  = 0.6853932584269663
$ws.Range("D9").Value  = 0.6460674157303371
$ws.Range("D15").Value = 0.6123595505617978
$ws.Range("D18").Value = 0.6235955056179775
$ws.Range("D21").Value = 0.6853932584269663
$ws.Range("D24").Value = 0.6910112359550562
$ws.Range("D27").Value = 0.5842696629213483
